# Updates cryptos list figures (price/volume) to the latest scrape,
# and fixes the Stellar/Toncoin row ordering (rows 24-25 swap content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.974.85"
$ws.Range("E2").Value = "  -1.28%  "
# Row 3
$ws.Range("D3").Value = "1.640.51"
$ws.Range("E3").Value = "  -0.54%  "
# Row 4
$ws.Range("E4").Value = "  +0.45%  "
# Row 5
$ws.Range("D5").Value = "'215.71"
$ws.Range("E5").Value = "  -0.92%  "
# Row 6
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -0.43%  "
# Row 7
$ws.Range("E7").Value = "  +0.36%  "
# Row 9
$ws.Range("E9").Value = "  -0.13%  "
# Row 10
$ws.Range("D10").Value = "'19.63"
$ws.Range("E10").Value = "  -1.74%  "
# Row 11
$ws.Range("E11").Value = "  +0.11%  "
# Row 12
$ws.Range("D12").Value = "1.866.62"
$ws.Range("E12").Value = "  -0.61%  "
# Row 13
$ws.Range("D13").Value = "'4.28"
$ws.Range("E13").Value = "  -0.77%  "
# Row 14
$ws.Range("D14").Value = "1.647.43"
$ws.Range("E14").Value = "  -0.64%  "
# Row 15
$ws.Range("E15").Value = "  -1.19%  "
# Row 16
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.55%  "
# Row 17
$ws.Range("D17").Value = "'63.04"
$ws.Range("E17").Value = "  -0.91%  "
# Row 18
$ws.Range("D18").Value = "25.931.75"
$ws.Range("E18").Value = "  -1.42%  "
# Row 19
$ws.Range("E19").Value = "  +0.39%  "
# Row 20
$ws.Range("D20").Value = "'193.09"
$ws.Range("E20").Value = "  -2.02%  "
# Row 21
$ws.Range("E21").Value = "  -1.91%  "
# Row 22
$ws.Range("E22").Value = "  -1.54%  "
# Row 24
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.131"
$ws.Range("E24").Value = "  +4.38%  "
# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  +0.27%  "
# Row 26
$ws.Range("D26").Value = "'143.60"
$ws.Range("E26").Value = "  +0.36%  "
# Row 27
$ws.Range("E27").Value = "  +0.42%  "
# Row 28
$ws.Range("D28").Value = "'6.87"
$ws.Range("E28").Value = "  -1.91%  "
# Row 29
$ws.Range("D29").Value = "'15.61"
$ws.Range("E29").Value = "  -0.45%  "
# Row 30
$ws.Range("E30").Value = "  -0.73%  "
# Row 31
$ws.Range("D31").Value = "'0.0504"
$ws.Range("E31").Value = "  -0.58%  "
# Row 32
$ws.Range("E32").Value = "  -1.74%  "
# Row 33
$ws.Range("E33").Value = "  -0.33%  "
# Row 35
$ws.Range("E35").Value = "  +1.77%  "
# Row 36
$ws.Range("D36").Value = "'0.902"
$ws.Range("E36").Value = "  -1.74%  "
# Row 37
$ws.Range("D37").Value = "1.135.81"
$ws.Range("E37").Value = "  -0.17%  "
# Row 38
$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  -2.14%  "
# Row 39
$ws.Range("D39").Value = "'2.47"
$ws.Range("E39").Value = "  -1.21%  "
# Row 40
$ws.Range("E40").Value = "  -0.25%  "
# Row 41
$ws.Range("E41").Value = "  +0.41%  "
# Row 42
$ws.Range("E42").Value = "  -3.09%  "
# Row 43
$ws.Range("D43").Value = "'99.42"
$ws.Range("E43").Value = "  -1.12%  "
# Row 44
$ws.Range("D44").Value = "'0.799"
$ws.Range("E44").Value = "  -0.59%  "
# Row 45
$ws.Range("D45").Value = "1.776.39"
$ws.Range("E45").Value = "  -0.65%  "
# Row 46
$ws.Range("E46").Value = "  +1.30%  "
# Row 47
$ws.Range("D47").Value = "'56.72"
$ws.Range("E47").Value = "  +0.35%  "
# Row 48
$ws.Range("E48").Value = "  +2.47%  "
# Row 49
$ws.Range("E49").Value = "  -2.09%  "
# Row 50
$ws.Range("D50").Value = "'7.69"
$ws.Range("E50").Value = "  -0.48%  "
# Row 51
$ws.Range("E51").Value = "  -0.54%  "

# The cells above that look like plain numbers ("215.71", "0.131", ...) were
# entered with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cells) instead of auto-converting to numeric values.
# Re-apply the Normal style so that text-forcing quote-prefix flag is cleared
# and formatting matches the untouched cells exactly.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
